$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 376.45456
$ws.Range("I19").Value = 366
$ws.Range("J19").Value = 398.85715
$ws.Range("K19").Value = 366
$ws.Range("L19").Value = 398.85715
$ws.Range("M19").Value = -191
$ws.Range("N19").Value = -748.85715

$ws.Range("H96").Value = 69665.125
$ws.Range("I96").Value = 3244.2727
$ws.Range("J96").Value = 215791
$ws.Range("K96").Value = 9732.8181
$ws.Range("L96").Value = 647373
$ws.Range("M96").Value = -8359.8181
$ws.Range("N96").Value = -650119

$ws.Range("H99").Value = 1010.35
$ws.Range("I99").Value = 826.17645
$ws.Range("K99").Value = 2478.52935
$ws.Range("M99").Value = -980.5293500000002

$ws.Range("H101").Value = 1104.2
$ws.Range("I101").Value = 875.5
$ws.Range("K101").Value = 2626.5
$ws.Range("M101").Value = -1004.5

$ws.Range("H133").Value = 59480
$ws.Range("J133").Value = 59480
$ws.Range("L133").Value = 59480
$ws.Range("N133").Value = -69600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12248.796
$ws.Range("I32").Value = 11935.889
$ws.Range("K32").Value = 11935.889
$ws.Range("M32").Value = -11648.889

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20540

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -21872

$ws.Range("H102").Value = 16774.666
$ws.Range("I102").Value = 2677.5
$ws.Range("K102").Value = 2677.5
$ws.Range("M102").Value = -1055.5

$ws.Range("H122").Value = 2102.611
$ws.Range("I122").Value = 2115.4375
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6346.3125
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3896.3125
$ws.Range("N122").Value = -10900

$ws.Range("H125").Value = 33230.145
$ws.Range("J125").Value = 33230.145
$ws.Range("L125").Value = 33230.145
$ws.Range("N125").Value = -43070.145

$ws.Range("H132").Value = 10639989
$ws.Range("I132").Value = 16130201
$ws.Range("J132").Value = 2704.625
$ws.Range("K132").Value = 48390603
$ws.Range("L132").Value = 8113.875
$ws.Range("M132").Value = -48388073
$ws.Range("N132").Value = -13173.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2626.7856
$ws.Range("I107").Value = 2396.2
$ws.Range("J107").Value = 3203.25
$ws.Range("K107").Value = 2396.2
$ws.Range("L107").Value = 3203.25
$ws.Range("M107").Value = -476.1999999999998
$ws.Range("N107").Value = -7043.25

$ws.Range("H122").Value = 40727.6
$ws.Range("J122").Value = 40727.6
$ws.Range("L122").Value = 40727.6
$ws.Range("N122").Value = -50527.6

$ws.Range("H124").Value = 50992
$ws.Range("J124").Value = 50992
$ws.Range("L124").Value = 50992
$ws.Range("N124").Value = -60812

$ws.Range("H132").Value = 31807.777
$ws.Range("J132").Value = 31807.777
$ws.Range("L132").Value = 31807.777
$ws.Range("N132").Value = -41927.777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4883.5566
$ws.Range("I31").Value = 2486.8333
$ws.Range("J31").Value = 5261.987
$ws.Range("K31").Value = 2486.8333
$ws.Range("L31").Value = 5261.987
$ws.Range("M31").Value = -2191.8333
$ws.Range("N31").Value = -5851.987

$ws.Range("H34").Value = 4883.5566
$ws.Range("I34").Value = 2486.8333
$ws.Range("J34").Value = 5261.987
$ws.Range("K34").Value = 2486.8333
$ws.Range("L34").Value = 5261.987
$ws.Range("M34").Value = -2284.8333
$ws.Range("N34").Value = -5665.987

$ws.Range("H52").Value = 57500
$ws.Range("J52").Value = 57500
$ws.Range("L52").Value = 57500
$ws.Range("N52").Value = -58088

$ws.Range("H112").Value = 45702
$ws.Range("J112").Value = 45702
$ws.Range("L112").Value = 45702
$ws.Range("N112").Value = -48656

$ws.Range("H137").Value = 56599.5
$ws.Range("J137").Value = 56599.5
$ws.Range("L137").Value = 56599.5
$ws.Range("N137").Value = -66799.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5443.3125
$ws.Range("I3").Value = 2714.3333
$ws.Range("J3").Value = 8952
$ws.Range("K3").Value = 8142.999899999999
$ws.Range("L3").Value = 26856
$ws.Range("M3").Value = -8030.999899999999
$ws.Range("N3").Value = -27080

$ws.Range("H23").Value = 432.69565
$ws.Range("I23").Value = 441.91666
$ws.Range("J23").Value = 422.63635
$ws.Range("K23").Value = 1325.74998
$ws.Range("L23").Value = 1267.90905
$ws.Range("M23").Value = -1090.74998
$ws.Range("N23").Value = -1737.90905

$ws.Range("H133").Value = 7466.6665
$ws.Range("I133").Value = 9100
$ws.Range("J133").Value = 7000
$ws.Range("K133").Value = 27300
$ws.Range("L133").Value = 21000
$ws.Range("M133").Value = -22240
$ws.Range("N133").Value = -31120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4467.222
$ws.Range("I97").Value = 2481
$ws.Range("J97").Value = 6950
$ws.Range("K97").Value = 2481
$ws.Range("L97").Value = 6950
$ws.Range("M97").Value = -1985
$ws.Range("N97").Value = -7942

$ws.Range("H122").Value = 1241
$ws.Range("J122").Value = 1340
$ws.Range("L122").Value = 4020
$ws.Range("N122").Value = -8920

$ws.Range("H124").Value = 38415
$ws.Range("J124").Value = 41768
$ws.Range("L124").Value = 41768
$ws.Range("N124").Value = -51588

$ws.Range("H135").Value = 49900
$ws.Range("J135").Value = 49900
$ws.Range("L135").Value = 49900
$ws.Range("N135").Value = -60040

$ws.Range("H136").Value = 85163
$ws.Range("J136").Value = 85163
$ws.Range("L136").Value = 255489
$ws.Range("N136").Value = -260589

$ws.Range("H138").Value = 54000
$ws.Range("J138").Value = 54000
$ws.Range("L138").Value = 54000
$ws.Range("N138").Value = -64280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 50975.145
$ws.Range("J134").Value = 50975.145
$ws.Range("L134").Value = 50975.145
$ws.Range("N134").Value = -61115.145

$ws.Range("H137").Value = 36648
$ws.Range("J137").Value = 36648
$ws.Range("L137").Value = 36648
$ws.Range("N137").Value = -46848

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 585.6
$ws.Range("I100").Value = 403.5
$ws.Range("J100").Value = 707
$ws.Range("K100").Value = 807
$ws.Range("L100").Value = 1414
$ws.Range("M100").Value = -266
$ws.Range("N100").Value = -2496

$ws.Range("H109").Value = 39373
$ws.Range("J109").Value = 39373
$ws.Range("L109").Value = 39373
$ws.Range("N109").Value = -42147

$ws.Range("H122").Value = 33036408
$ws.Range("I122").Value = 33036408
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 99109224
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -99106774
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 45999.5
$ws.Range("J137").Value = 45999.5
$ws.Range("L137").Value = 45999.5
$ws.Range("N137").Value = -56199.5

$ws.Range("H139").Value = 56350
$ws.Range("J139").Value = 56350
$ws.Range("L139").Value = 56350
$ws.Range("N139").Value = -66630
